$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on D and E columns for the affected rows so that numeric-looking
# strings (prices like "29.368.38", "1.000", "0.00000000117") and percentage strings
# retain their exact original text representation instead of being parsed as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '29.368.38'
$ws.Range("E2").Value = '  -0.44%  '
$ws.Range("D3").Value = '1.846.52'
$ws.Range("E3").Value = '  -0.22%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '240.28'
$ws.Range("E5").Value = '  -0.64%  '
$ws.Range("D6").Value = '0.6302'
$ws.Range("E6").Value = '  +0.13%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("B8").Value = 'OKB'
$ws.Range("C8").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D8").Value = '47.79'
$ws.Range("E8").Value = '  -0.70%  '
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").Value = '0.07537'
$ws.Range("E9").Value = '  +0.06%  '
$ws.Range("B10").Value = 'Cardano'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D10").Value = '0.2955'
$ws.Range("E10").Value = '  -0.75%  '
$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").Value = '24.41'
$ws.Range("E11").Value = '  +0.04%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = '0.07702'
$ws.Range("E12").Value = '  -0.29%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.868.09'
$ws.Range("E13").Value = '  +0.78%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '4.993'
$ws.Range("E14").Value = '  -0.19%  '
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").Value = '0.6831'
$ws.Range("E15").Value = '  -1.38%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = '0.000009998'
$ws.Range("E16").Value = '  +2.26%  '
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").Value = '82.86'
$ws.Range("E17").Value = '  -0.85%  '
$ws.Range("D18").Value = '6.133'
$ws.Range("E18").Value = '  -1.85%  '
$ws.Range("D19").Value = '29.398.00'
$ws.Range("E19").Value = '  -0.46%  '
$ws.Range("D20").Value = '227.63'
$ws.Range("E20").Value = '  -2.47%  '
$ws.Range("E21").Value = '  -0.34%  '
$ws.Range("E22").Value = '  +0.03%  '
$ws.Range("D23").Value = '7.551'
$ws.Range("E23").Value = '  -1.17%  '
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("D25").Value = '157.36'
$ws.Range("E25").Value = '  +1.69%  '
$ws.Range("E26").Value = '  +0.42%  '
$ws.Range("D27").Value = '8.365'
$ws.Range("E27").Value = '  -1.11%  '
$ws.Range("E28").Value = '  -0.35%  '
$ws.Range("D29").Value = '1.465'
$ws.Range("E29").Value = '  -0.83%  '
$ws.Range("D30").Value = '1.262'
$ws.Range("E30").Value = '  +0.89%  '
$ws.Range("D31").Value = '0.05684'
$ws.Range("E31").Value = '  -4.37%  '
$ws.Range("E32").Value = '  +0.17%  '
$ws.Range("D33").Value = '4.019'
$ws.Range("E33").Value = '  -0.39%  '
$ws.Range("D34").Value = '1.845'
$ws.Range("E34").Value = '  -2.35%  '
$ws.Range("E35").Value = '  -1.21%  '
$ws.Range("D36").Value = '0.7123'
$ws.Range("E36").Value = '  -0.97%  '
$ws.Range("D37").Value = '2.593'
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("D38").Value = '1.261.50'
$ws.Range("E38").Value = '  +1.81%  '
$ws.Range("D39").Value = '0.01814'
$ws.Range("E39").Value = '  +0.92%  '
$ws.Range("D40").Value = '2.773'
$ws.Range("E40").Value = '  -0.80%  '
$ws.Range("D41").Value = '0.9063'
$ws.Range("E41").Value = '  -0.13%  '
$ws.Range("D42").Value = '6.164'
$ws.Range("E42").Value = '  +0.41%  '
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("D44").Value = '101.34'
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("E45").Value = '  -1.39%  '
$ws.Range("D46").Value = '7.072'
$ws.Range("E46").Value = '  -4.41%  '
$ws.Range("D47").Value = '0.4042'
$ws.Range("E47").Value = '  -0.09%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '9.058'
$ws.Range("E48").Value = '  -1.26%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = '1.684'
$ws.Range("E49").Value = '  -1.00%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = '0.1120'
$ws.Range("E50").Value = '  +0.10%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.05734'
$ws.Range("E51").Value = '  -0.44%  '
